# Fruta / hortaliza, semanal
# Insert two new weekly data rows into the "Zapallo italiano" (Vega Modelo de
# Temuco) price sheet. Inserting whole rows at 618:619 shifts all the
# existing data (old rows 618-702) down by two, which matches the target
# dimension growing from A1:R702 to A1:R704, and carries the date-number
# formatting (style index on column D) down with the shifted rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("618:619").Insert()

# New row 618
$ws.Range("A618").Value = 10
$ws.Range("B618").Value = "Vega Modelo de Temuco"
$ws.Range("C618").Value = "La Araucanía"
$ws.Range("D618").Value = 44984
$ws.Range("E618").Value = 9
$ws.Range("F618").Value = 100112032
$ws.Range("G618").Value = "Zapallo italiano"
$ws.Range("H618").Value = "Sin especificar"
$ws.Range("I618").Value = "Primera"
$ws.Range("J618").Value = 50
$ws.Range("K618").Value = 9000
$ws.Range("L618").Value = 9000
$ws.Range("M618").Value = 9000
$ws.Range("N618").Value = "$/caja 36 unidades"
$ws.Range("O618").Value = "Región de La Araucanía"
$ws.Range("P618").Value = 250
$ws.Range("Q618").Value = 36
$ws.Range("R618").Value = "Hortaliza"

# New row 619
$ws.Range("A619").Value = 10
$ws.Range("B619").Value = "Vega Modelo de Temuco"
$ws.Range("C619").Value = "La Araucanía"
$ws.Range("D619").Value = 44984
$ws.Range("E619").Value = 9
$ws.Range("F619").Value = 100112032
$ws.Range("G619").Value = "Zapallo italiano"
$ws.Range("H619").Value = "Sin especificar"
$ws.Range("I619").Value = "Primera"
$ws.Range("J619").Value = 240
$ws.Range("K619").Value = 10000
$ws.Range("L619").Value = 12000
$ws.Range("M619").Value = 11000
$ws.Range("N619").Value = "$/caja 50 unidades"
$ws.Range("O619").Value = "Región del Maule"
$ws.Range("P619").Value = 220
$ws.Range("Q619").Value = 50
$ws.Range("R619").Value = "Hortaliza"
